$d = $word.ActiveDocument

# (1) Tidy the stray leading/trailing spaces around "Phonemic" in the
#     --glm_fea " Phonemic " argument of the step2_time_cluster / Resp_inRep
#     command line, turning it into --glm_fea "Phonemic".
$d.Content.Find.Execute(" Phonemic ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Phonemic", 2)

# (2) Strip the stray paragraph-mark formatting (an eastAsia rFonts hint
#     left on an otherwise empty paragraph's mark) so the paragraph becomes
#     a plain empty paragraph like its neighbours. Locate it by inspecting
#     each empty paragraph's own OOXML for that exact pPr/rPr/rFonts hint.
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.Length -gt 1) { continue }
    $xml = $p.Range.WordOpenXML
    $m = [regex]::Match($xml, '<w:body>(.*?)<w:sectPr')
    if (-not $m.Success) { continue }
    $body = $m.Groups[1].Value
    if ($body -like '*<w:pPr><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr></w:pPr>*') {
        $p.Range.Select()
        $word.Selection.ClearFormatting()
        break
    }
}
